# ExtractCVs.xaml was amended so the %CVMatch variable is built as an
# integer instead of a double -> the matchCVPercentage_col column (D) now
# holds whole numbers (rounded), not 2-decimal percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# matchCVPercentage_col (column D) - values rounded to the nearest integer
$percentages = @{
    2  = 41
    3  = 41
    4  = 35
    5  = 35
    6  = 30
    7  = 19
    8  = 14
    9  = 14
    10 = 3
    11 = 3
    12 = 3
    13 = 3
    14 = 3
    15 = 3
}

foreach ($row in $percentages.Keys) {
    $ws.Cells.Item($row, 4).Value = $percentages[$row]
}

# Restore the workbook window to a maximized/full-size layout and scroll /
# reselect as it was left after the edit (view-state bookkeeping).
$excel.ActiveWindow.WindowState = -4137  # xlMaximized
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A10").Select()
